# Fixed update to excel issue
# 1. Rename "Requested quantity" headers to dataset-specific names
# 2. Add a new "PO Forecast" worksheet with forecast data (ds, PO_Forecast, yhat_lower, yhat_upper)

$wb = $excel.ActiveWorkbook

# --- Step 1: rename headers on existing sheets -----------------------------
$wsWeekly = $wb.Worksheets.Item(1)   # "Weekly Quantity"
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item(2)  # "Monthly Trend"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Step 2: create the new "PO Forecast" worksheet at the end -------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Forecast data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$forecastData = @(
  @(44941.99999999999, 154, 66.64667518387783, 245.2917163484848),
  @(44990.99999999999, 113, 23.66901044722736, 206.0355957847412),
  @(44997.99999999999, 108, 18.1646399270352, 196.0562441034052),
  @(45004.99999999999, 102, 8.575117671241163, 196.8211233126173),
  @(45011.99999999999, 96, 6.718584329734803, 182.8771650578433),
  @(45018.99999999999, 90, 2.013426710704346, 177.8848702301342),
  @(45025.99999999999, 84, -8.199653316591911, 170.4326461896836),
  @(45032.99999999999, 79, -13.2508490680484, 164.1709991846542),
  @(45046.99999999999, 67, -25.46044564532125, 156.8764118863134),
  @(45053.99999999999, 61, -24.59916978179405, 150.6234302824076),
  @(45060.99999999999, 55, -37.29497393621009, 149.5496867032041),
  @(45067.99999999999, 50, -41.68327971620437, 137.72042972618),
  @(45074.99999999999, 44, -40.53864606650357, 131.0141990782542),
  @(45081.99999999999, 38, -47.39856438162433, 126.3924545461688),
  @(45109.99999999999, 15, -75.56567095586924, 106.1266264402587),
  @(45116.99999999999, 9, -78.89004383049351, 100.7535135411421),
  @(45123.99999999999, 3, -82.90860802206824, 97.00222412033571),
  @(45130.99999999999, 0, -89.29047862616092, 83.47932226721601),
  @(45137.99999999999, 0, -96.64985894986765, 82.18657855625456),
  @(45144.99999999999, 0, -104.5604476617452, 71.18954415786791),
  @(45151.99999999999, 0, -103.6408580790045, 60.75939136988941),
  @(45158.99999999999, 0, -116.4084572110623, 62.27176552684252),
  @(45165.99999999999, 0, -119.2862881485677, 55.33959952294895),
  @(45172.99999999999, 0, -121.6553053385115, 52.33076441433511),
  @(45179.99999999999, 0, -130.0593672999628, 47.8453093169168),
  @(45186.99999999999, 0, -132.164936462561, 36.90865804013835),
  @(45193.99999999999, 0, -137.988526149477, 38.594261379259),
  @(45200.99999999999, 0, -149.1012826804936, 29.50437135707439)
)

$rowIndex = 2
foreach ($dataRow in $forecastData) {
    $wsForecast.Cells.Item($rowIndex, 1).Value = $dataRow[0]
    $wsForecast.Cells.Item($rowIndex, 2).Value = $dataRow[1]
    $wsForecast.Cells.Item($rowIndex, 3).Value = $dataRow[2]
    $wsForecast.Cells.Item($rowIndex, 4).Value = $dataRow[3]
    $rowIndex++
}

# --- Step 3: replicate formatting from the "Weekly Quantity" sheet ---------
# Header row formatting (bold, bordered, centered) -> columns A:D
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

# Date-column formatting (numFmt yyyy-mm-dd hh:mm:ss) -> column A, rows 2:29
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A29").PasteSpecial(-4122)  # xlPasteFormats

$wsForecast.Range("A1").Select()
$excel.CutCopyMode = $false
